$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3 into row 4 (values + styles + row height), matching the
# template's pattern for new user rows.
$ws.Range("A3:P3").Copy($ws.Range("A4:P4"))
$ws.Rows.Item(4).RowHeight = 16

# Update the e-mail cell (C4) to a new address and turn it into a mailto
# hyperlink, restoring the original (non-hyperlink-styled) cell formatting
# afterwards so it matches the template's existing hyperlink cells.
$ws.Range("C4").Value = "user_id_3@1.c"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:user_id_3@1.c")
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# The password cell (D4) keeps the same value as D3, but still needs its
# own hyperlink relationship/entry (mirroring D2/D3).
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:Calong@2015")
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# Move the active selection, matching the saved cursor position.
$ws.Range("C18").Select()
